$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Complete Custom Extent values in column B (rows 2-15) from 10.5 to 11.5
$ws.Range("B2:B15").Value = 11.5

# Update the view: scroll back to top (clear topLeftCell) and select B2:B15
$ws.Range("B2:B15").Select()
